$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.911.05'
$ws.Range('E2').Value = '  +0.85%  '
$ws.Range('D3').Value = '3.509.03'
$ws.Range('E3').Value = '  -0.42%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '600.60'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +0.55%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '182.77'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +5.22%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.597'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +0.40%  '
$ws.Range('E9').Value = '  +5.79%  '
$ws.Range('E10').Value = '  -2.19%  '
$ws.Range('E11').Value = '  -0.46%  '
$ws.Range('D12').Value = '4.119.01'
$ws.Range('E12').Value = '  -0.39%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '32.56'
$ws.Range('D13').ClearFormats()
$ws.Range('E14').Value = '  -0.17%  '
$ws.Range('E15').Value = '  +0.16%  '
$ws.Range('D16').Value = '67.887.04'
$ws.Range('E16').Value = '  +0.87%  '
$ws.Range('D17').Value = '3.522.66'
$ws.Range('E17').Value = '  -0.02%  '
$ws.Range('E18').Value = '  +0.82%  '
$ws.Range('E19').Value = '  +3.04%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '397.19'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -0.02%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '8.09'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +1.20%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '73.70'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +0.37%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.545'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +1.01%  '
$ws.Range('E24').Value = '  +0.16%  '
$ws.Range('E25').Value = '  +0.05%  '
$ws.Range('E26').Value = '  +1.00%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.40'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +1.07%  '
$ws.Range('E28').Value = '  -0.69%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.01'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +1.45%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.31'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +0.21%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.46'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -0.09%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.07'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -0.64%  '
$ws.Range('E33').Value = '  -0.67%  '
$ws.Range('E34').Value = '  +0.15%  '
$ws.Range('E35').Value = '  +0.07%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.68'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +1.74%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '163.53'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +0.14%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.96'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +2.24%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.876'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -2.25%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '7.12'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +2.42%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.75'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +0.63%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '27.83'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +0.73%  '
$ws.Range('E43').Value = '  +2.78%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '26.69'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +0.79%  '
$ws.Range('E45').Value = '  -1.35%  '
$ws.Range('D46').Value = '2.815.24'
$ws.Range('E46').Value = '  +0.52%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '42.42'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -1.15%  '
$ws.Range('E48').Value = '  -0.94%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '345.55'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +1.16%  '
$ws.Range('E50').Value = '  -1.14%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '33.78'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +0.09%  '
